$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set manually-assigned category values (column B) for rows that were categorized ---
$ws.Range("B11").Value = "RUIM"
$ws.Range("B12").Value = "RUIM"
$ws.Range("B13").Value = "BOM"
$ws.Range("B14").Value = "RUIM"
$ws.Range("B15").Value = "BOM"
$ws.Range("B16").Value = "CUPOM"
$ws.Range("B17").Value = "RUIM"
$ws.Range("B18").Value = "BOM"
$ws.Range("B19").Value = "BOM"
$ws.Range("B20").Value = "BOM"
$ws.Range("B21").Value = "RUIM"
$ws.Range("B22").Value = "BOM"
$ws.Range("B23").Value = "CUPOM"
$ws.Range("B24").Value = "BOM"
$ws.Range("B25").Value = "IRRELEVANTE"
$ws.Range("B26").Value = "RUIM"
$ws.Range("B27").Value = "CUPOM"
$ws.Range("B28").Value = "BOM"
$ws.Range("B29").Value = "IRRELEVANTE"
$ws.Range("B30").Value = "RUIM"
$ws.Range("B31").Value = "BOM"
$ws.Range("B32").Value = "BOM"
$ws.Range("B33").Value = "BOM"
$ws.Range("B34").Value = "CUPOM"
$ws.Range("B35").Value = "BOM"
$ws.Range("B36").Value = "BOM"
$ws.Range("B37").Value = "BOM"
$ws.Range("B38").Value = "IRRELEVANTE"
$ws.Range("B39").Value = "BOM"
$ws.Range("B40").Value = "BOM"
$ws.Range("B41").Value = "CUPOM"
$ws.Range("B42").Value = "CUPOM"
$ws.Range("B43").Value = "IRRELEVANTE"
$ws.Range("B44").Value = "RUIM"
$ws.Range("B45").Value = "RUIM"
$ws.Range("B46").Value = "BOM"
$ws.Range("B47").Value = "BOM"
$ws.Range("B48").Value = "BOM"
$ws.Range("B49").Value = "BOM"
$ws.Range("B50").Value = "RUIM"
$ws.Range("B51").Value = "BOM"
$ws.Range("B52").Value = "RUIM"
$ws.Range("B53").Value = "IRRELEVANTE"
$ws.Range("B233").Value = "CUPOM"

# --- Row 44: wrap text + taller row (reformatted entry) ---
$ws.Range("A44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 29

# --- Column width adjustments ---
$ws.Columns.Item(1).ColumnWidth = 167.83333333333334
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 33

# --- View / window state ---
$excel.Goto($ws.Range("A11"), $false)
$ws.Range("B54").Select()
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.WindowState = -4143
